$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the full data range (rank/index/name/sem1/sem2/sem3/cgpa) ---
# Several students' sem3 & CGPA values were corrected (EDP + NM recorrection),
# which changes the CGPA-descending sort order of the whole table.
$data = New-Object 'object[,]' 115,7
$data[0,0] = 1
$data[0,1] = 230266
$data[0,2] = 'JATHUNWATHTHA J.C.R.N.'
$data[0,3] = 4
$data[0,4] = 4
$data[0,5] = 4
$data[0,6] = 4
$data[1,0] = 2
$data[1,1] = 230487
$data[1,2] = 'PERERA W.A.L.S.'
$data[1,3] = 4
$data[1,4] = 4
$data[1,5] = 4
$data[1,6] = 4
$data[2,0] = 3
$data[2,1] = 230018
$data[2,2] = 'ADIKARAM D.M.G.H.'
$data[2,3] = 4
$data[2,4] = 4
$data[2,5] = 3.96
$data[2,6] = 3.99
$data[3,0] = 4
$data[3,1] = 230074
$data[3,2] = 'BANDARA H.Y.W.'
$data[3,3] = 4
$data[3,4] = 4
$data[3,5] = 3.96
$data[3,6] = 3.99
$data[4,0] = 5
$data[4,1] = 230082
$data[4,2] = 'BANDARA W.D.A.C.'
$data[4,3] = 4
$data[4,4] = 4
$data[4,5] = 3.96
$data[4,6] = 3.99
$data[5,0] = 6
$data[5,1] = 230171
$data[5,2] = 'ELAPATHA C.D.'
$data[5,3] = 4
$data[5,4] = 4
$data[5,5] = 3.97
$data[5,6] = 3.99
$data[6,0] = 7
$data[6,1] = 230436
$data[6,2] = 'NETTIKUMARA N.A.H.G.'
$data[6,3] = 4
$data[6,4] = 3.96
$data[6,5] = 4
$data[6,6] = 3.99
$data[7,0] = 8
$data[7,1] = 230476
$data[7,2] = 'PERERA G.M.B.'
$data[7,3] = 4
$data[7,4] = 4
$data[7,5] = 3.96
$data[7,6] = 3.99
$data[8,0] = 9
$data[8,1] = 230138
$data[8,2] = 'DHANANJAYA K.T.G.T.N.'
$data[8,3] = 4
$data[8,4] = 4
$data[8,5] = 3.96
$data[8,6] = 3.98
$data[9,0] = 10
$data[9,1] = 230548
$data[9,2] = 'RATNAYAKE R.M.S.H.'
$data[9,3] = 4
$data[9,4] = 4
$data[9,5] = 3.96
$data[9,6] = 3.98
$data[10,0] = 11
$data[10,1] = 230689
$data[10,2] = 'WEERAKOON A.H.T.M.'
$data[10,3] = 4
$data[10,4] = 4
$data[10,5] = 3.95
$data[10,6] = 3.98
$data[11,0] = 12
$data[11,1] = 230108
$data[11,2] = 'COLOMBAGE D.M.'
$data[11,3] = 3.94
$data[11,4] = 4
$data[11,5] = 3.95
$data[11,6] = 3.97
$data[12,0] = 13
$data[12,1] = 230256
$data[12,2] = 'ILANKOON I.M.M.K.B.'
$data[12,3] = 4
$data[12,4] = 4
$data[12,5] = 3.93
$data[12,6] = 3.97
$data[13,0] = 14
$data[13,1] = 230318
$data[13,2] = 'KARIYAWASAM J.H.D.'
$data[13,3] = 4
$data[13,4] = 4
$data[13,5] = 3.92
$data[13,6] = 3.97
$data[14,0] = 15
$data[14,1] = 230352
$data[14,2] = 'KUMARA K.B.R.'
$data[14,3] = 3.94
$data[14,4] = 4
$data[14,5] = 3.96
$data[14,6] = 3.97
$data[15,0] = 16
$data[15,1] = 230481
$data[15,2] = 'PERERA K.W.A.O.V.'
$data[15,3] = 3.85
$data[15,4] = 4
$data[15,5] = 4
$data[15,6] = 3.97
$data[16,0] = 17
$data[16,1] = 230544
$data[16,2] = 'RATHNAYAKE M.A.G.K.N.'
$data[16,3] = 4
$data[16,4] = 4
$data[16,5] = 3.92
$data[16,6] = 3.97
$data[17,0] = 18
$data[17,1] = 230121
$data[17,2] = 'DE MEL D.J.'
$data[17,3] = 3.96
$data[17,4] = 4
$data[17,5] = 3.92
$data[17,6] = 3.96
$data[18,0] = 19
$data[18,1] = 230355
$data[18,2] = 'KUMARASINGHE M.N.'
$data[18,3] = 4
$data[18,4] = 4
$data[18,5] = 3.9
$data[18,6] = 3.96
$data[19,0] = 20
$data[19,1] = 230469
$data[19,2] = 'PEIRIS E.A.S.S.'
$data[19,3] = 4
$data[19,4] = 3.96
$data[19,5] = 3.94
$data[19,6] = 3.96
$data[20,0] = 21
$data[20,1] = 230680
$data[20,2] = 'WANIGASUNDARA W.M.H.'
$data[20,3] = 4
$data[20,4] = 4
$data[20,5] = 3.9
$data[20,6] = 3.96
$data[21,0] = 22
$data[21,1] = 230155
$data[21,2] = 'DISSANAYAKA D.M.D.P.'
$data[21,3] = 4
$data[21,4] = 3.93
$data[21,5] = 3.96
$data[21,6] = 3.95
$data[22,0] = 23
$data[22,1] = 230159
$data[22,2] = 'DISSANAYAKE G.R.G.K.'
$data[22,3] = 4
$data[22,4] = 3.96
$data[22,5] = 3.9
$data[22,6] = 3.95
$data[23,0] = 24
$data[23,1] = 230186
$data[23,2] = 'FERNANDO W.H.D.'
$data[23,3] = 4
$data[23,4] = 4
$data[23,5] = 3.87
$data[23,6] = 3.95
$data[24,0] = 25
$data[24,1] = 230390
$data[24,2] = 'MALDENIYA P.A.D.G.R.'
$data[24,3] = 4
$data[24,4] = 4
$data[24,5] = 3.87
$data[24,6] = 3.95
$data[25,0] = 26
$data[25,1] = 230508
$data[25,2] = 'RAHUL B.'
$data[25,3] = 4
$data[25,4] = 4
$data[25,5] = 3.87
$data[25,6] = 3.95
$data[26,0] = 27
$data[26,1] = 230051
$data[26,2] = 'ARACHCHI A.D.I.D.'
$data[26,3] = 4
$data[26,4] = 4
$data[26,5] = 3.83
$data[26,6] = 3.94
$data[27,0] = 28
$data[27,1] = 230258
$data[27,2] = 'IMADUWAGE O.N.H.'
$data[27,3] = 3.89
$data[27,4] = 3.96
$data[27,5] = 3.96
$data[27,6] = 3.94
$data[28,0] = 29
$data[28,1] = 230468
$data[28,2] = 'PATHIRANA P.T.S.'
$data[28,3] = 3.91
$data[28,4] = 4
$data[28,5] = 3.9
$data[28,6] = 3.94
$data[29,0] = 30
$data[29,1] = 230566
$data[29,2] = 'SAMARASINGHE S.M.R.R.'
$data[29,3] = 3.96
$data[29,4] = 4
$data[29,5] = 3.88
$data[29,6] = 3.94
$data[30,0] = 31
$data[30,1] = 230332
$data[30,2] = 'KEERAWELLA K.P.C.P.'
$data[30,3] = 4
$data[30,4] = 4
$data[30,5] = 3.79
$data[30,6] = 3.93
$data[31,0] = 32
$data[31,1] = 230486
$data[31,2] = 'PERERA U.I.H.'
$data[31,3] = 4
$data[31,4] = 4
$data[31,5] = 3.83
$data[31,6] = 3.93
$data[32,0] = 33
$data[32,1] = 230140
$data[32,2] = 'DHARMAKEERTHI P.K.G.C.L.'
$data[32,3] = 3.94
$data[32,4] = 3.96
$data[32,5] = 3.87
$data[32,6] = 3.92
$data[33,0] = 34
$data[33,1] = 230521
$data[33,2] = 'RANASINGHE D.P.H.'
$data[33,3] = 4
$data[33,4] = 4
$data[33,5] = 3.79
$data[33,6] = 3.92
$data[34,0] = 35
$data[34,1] = 230536
$data[34,2] = 'RASANJANA W.P.G.R.A.'
$data[34,3] = 3.96
$data[34,4] = 3.96
$data[34,5] = 3.85
$data[34,6] = 3.92
$data[35,0] = 36
$data[35,1] = 230197
$data[35,2] = 'GARUSINGHE S.B.'
$data[35,3] = 4
$data[35,4] = 3.88
$data[35,5] = 3.9
$data[35,6] = 3.91
$data[36,0] = 37
$data[36,1] = 230322
$data[36,2] = 'KARUNARATHNA G.K.T.'
$data[36,3] = 4
$data[36,4] = 3.92
$data[36,5] = 3.86
$data[36,6] = 3.91
$data[37,0] = 38
$data[37,1] = 230687
$data[37,2] = 'WEDAMESTRIGE A.N.'
$data[37,3] = 4
$data[37,4] = 3.93
$data[37,5] = 3.83
$data[37,6] = 3.91
$data[38,0] = 39
$data[38,1] = 230100
$data[38,2] = 'CHANDRAKUMARA H.A.D.C.'
$data[38,3] = 4
$data[38,4] = 4
$data[38,5] = 3.72
$data[38,6] = 3.9
$data[39,0] = 40
$data[39,1] = 230724
$data[39,2] = 'WIJESEKARA W.A.G.S.'
$data[39,3] = 4
$data[39,4] = 3.9
$data[39,5] = 3.83
$data[39,6] = 3.9
$data[40,0] = 41
$data[40,1] = 230038
$data[40,2] = 'AMARATHUNGE A.M.N.L.'
$data[40,3] = 4
$data[40,4] = 4
$data[40,5] = 3.69
$data[40,6] = 3.89
$data[41,0] = 42
$data[41,1] = 230417
$data[41,2] = 'MUNASINGHE A.I.'
$data[41,3] = 4
$data[41,4] = 3.92
$data[41,5] = 3.8
$data[41,6] = 3.89
$data[42,0] = 43
$data[42,1] = 230130
$data[42,2] = 'DESHAN W.U.'
$data[42,3] = 4
$data[42,4] = 3.96
$data[42,5] = 3.71
$data[42,6] = 3.88
$data[43,0] = 44
$data[43,1] = 230145
$data[43,2] = 'DILHAN W.A.'
$data[43,3] = 3.94
$data[43,4] = 4
$data[43,5] = 3.7
$data[43,6] = 3.88
$data[44,0] = 45
$data[44,1] = 230212
$data[44,2] = 'GUNASEKARA L.U.A.'
$data[44,3] = 3.96
$data[44,4] = 3.88
$data[44,5] = 3.82
$data[44,6] = 3.88
$data[45,0] = 46
$data[45,1] = 230300
$data[45,2] = 'JAYAWEERA N.S.'
$data[45,3] = 3.94
$data[45,4] = 4
$data[45,5] = 3.71
$data[45,6] = 3.88
$data[46,0] = 47
$data[46,1] = 230321
$data[46,2] = 'KARUNANAYAKE A.H.D.'
$data[46,3] = 4
$data[46,4] = 3.95
$data[46,5] = 3.73
$data[46,6] = 3.88
$data[47,0] = 48
$data[47,1] = 230477
$data[47,2] = 'PERERA H.A.J.I.'
$data[47,3] = 3.94
$data[47,4] = 3.83
$data[47,5] = 3.87
$data[47,6] = 3.87
$data[48,0] = 49
$data[48,1] = 230659
$data[48,2] = 'UPEKSHANI T.S.'
$data[48,3] = 3.86
$data[48,4] = 3.96
$data[48,5] = 3.8
$data[48,6] = 3.87
$data[49,0] = 50
$data[49,1] = 230697
$data[49,2] = 'WEERASINGHE J.A.H.R.'
$data[49,3] = 3.96
$data[49,4] = 4
$data[49,5] = 3.67
$data[49,6] = 3.87
$data[50,0] = 51
$data[50,1] = 230045
$data[50,2] = 'ANTHONY C.S.B.'
$data[50,3] = 4
$data[50,4] = 3.92
$data[50,5] = 3.7
$data[50,6] = 3.86
$data[51,0] = 52
$data[51,1] = 230058
$data[51,2] = 'AROSHANA H.A.P.'
$data[51,3] = 4
$data[51,4] = 3.92
$data[51,5] = 3.71
$data[51,6] = 3.86
$data[52,0] = 53
$data[52,1] = 230065
$data[52,2] = 'AYANAJA N.B.G.M.'
$data[52,3] = 3.89
$data[52,4] = 3.84
$data[52,5] = 3.85
$data[52,6] = 3.86
$data[53,0] = 54
$data[53,1] = 230613
$data[53,2] = 'SHEHAN M.N.N.'
$data[53,3] = 4
$data[53,4] = 3.95
$data[53,5] = 3.67
$data[53,6] = 3.86
$data[54,0] = 55
$data[54,1] = 230211
$data[54,2] = 'GUNASEKARA K.S.'
$data[54,3] = 4
$data[54,4] = 3.89
$data[54,5] = 3.72
$data[54,6] = 3.85
$data[55,0] = 56
$data[55,1] = 230539
$data[55,2] = 'RATHEESHAN A.R.'
$data[55,3] = 3.94
$data[55,4] = 4
$data[55,5] = 3.65
$data[55,6] = 3.85
$data[56,0] = 57
$data[56,1] = 230492
$data[56,2] = 'PITIWADUGE D.N.'
$data[56,3] = 3.94
$data[56,4] = 3.92
$data[56,5] = 3.69
$data[56,6] = 3.84
$data[57,0] = 58
$data[57,1] = 230500
$data[57,2] = 'PRISHMIKA H.W.N.'
$data[57,3] = 3.96
$data[57,4] = 3.9
$data[57,5] = 3.68
$data[57,6] = 3.83
$data[58,0] = 59
$data[58,1] = 230629
$data[58,2] = 'TENNAKOON U.G.R.B.'
$data[58,3] = 3.96
$data[58,4] = 3.9
$data[58,5] = 3.67
$data[58,6] = 3.83
$data[59,0] = 60
$data[59,1] = 230180
$data[59,2] = 'FERNANDO H.M.D.'
$data[59,3] = 3.94
$data[59,4] = 3.83
$data[59,5] = 3.74
$data[59,6] = 3.82
$data[60,0] = 61
$data[60,1] = 230353
$data[60,2] = 'KUMARA P.K.M.P.'
$data[60,3] = 3.9
$data[60,4] = 3.96
$data[60,5] = 3.61
$data[60,6] = 3.82
$data[61,0] = 62
$data[61,1] = 230470
$data[61,2] = 'PEIRIS T.S.R.'
$data[61,3] = 4
$data[61,4] = 4
$data[61,5] = 3.53
$data[61,6] = 3.82
$data[62,0] = 63
$data[62,1] = 230502
$data[62,2] = 'PRIYADARSHANA S.A.D.'
$data[62,3] = 4
$data[62,4] = 3.95
$data[62,5] = 3.59
$data[62,6] = 3.82
$data[63,0] = 64
$data[63,1] = 230218
$data[63,2] = 'GUNATHUNGA U.A.'
$data[63,3] = 3.94
$data[63,4] = 3.77
$data[63,5] = 3.78
$data[63,6] = 3.81
$data[64,0] = 65
$data[64,1] = 230525
$data[64,2] = 'RANATHUNGA R.J.K.O.H.'
$data[64,3] = 4
$data[64,4] = 4
$data[64,5] = 3.5
$data[64,6] = 3.81
$data[65,0] = 66
$data[65,1] = 230726
$data[65,2] = 'WIJESINGHE U.G.S.K.D.'
$data[65,3] = 3.89
$data[65,4] = 3.9
$data[65,5] = 3.66
$data[65,6] = 3.81
$data[66,0] = 67
$data[66,1] = 230063
$data[66,2] = 'ATHUKORALA U.R.'
$data[66,3] = 4
$data[66,4] = 3.92
$data[66,5] = 3.56
$data[66,6] = 3.8
$data[67,0] = 68
$data[67,1] = 230070
$data[67,2] = 'BALASOORIYA B.R.B.D.'
$data[67,3] = 3.96
$data[67,4] = 3.82
$data[67,5] = 3.67
$data[67,6] = 3.8
$data[68,0] = 69
$data[68,1] = 230016
$data[68,2] = 'ABISHEK L.'
$data[68,3] = 4
$data[68,4] = 3.91
$data[68,5] = 3.48
$data[68,6] = 3.79
$data[69,0] = 70
$data[69,1] = 230164
$data[69,2] = 'DISSANAYAKE R.K.T.'
$data[69,3] = 3.96
$data[69,4] = 3.88
$data[69,5] = 3.59
$data[69,6] = 3.79
$data[70,0] = 71
$data[70,1] = 230195
$data[70,2] = 'GAMAGE SK'
$data[70,3] = 3.96
$data[70,4] = 3.8
$data[70,5] = 3.68
$data[70,6] = 3.79
$data[71,0] = 72
$data[71,1] = 230280
$data[71,2] = 'JAYASINGHE J.A.P.R.'
$data[71,3] = 3.85
$data[71,4] = 3.92
$data[71,5] = 3.6
$data[71,6] = 3.79
$data[72,0] = 73
$data[72,1] = 230507
$data[72,2] = 'RAHMAN M.F.A.'
$data[72,3] = 3.86
$data[72,4] = 3.87
$data[72,5] = 3.65
$data[72,6] = 3.79
$data[73,0] = 74
$data[73,1] = 230585
$data[73,2] = 'SARUKA U.'
$data[73,3] = 3.94
$data[73,4] = 3.96
$data[73,5] = 3.54
$data[73,6] = 3.79
$data[74,0] = 75
$data[74,1] = 230526
$data[74,2] = 'RANAWAKA R.A.C.D.'
$data[74,3] = 3.94
$data[74,4] = 3.95
$data[74,5] = 3.5
$data[74,6] = 3.77
$data[75,0] = 76
$data[75,1] = 230020
$data[75,2] = 'AHAMED A.M.S.'
$data[75,3] = 4
$data[75,4] = 3.87
$data[75,5] = 3.46
$data[75,6] = 3.76
$data[76,0] = 77
$data[76,1] = 230052
$data[76,2] = 'ARACHCHIGE M. A. D. T. S.'
$data[76,3] = 3.75
$data[76,4] = 3.85
$data[76,5] = 3.63
$data[76,6] = 3.76
$data[77,0] = 78
$data[77,1] = 230654
$data[77,2] = 'UMAIR A.'
$data[77,3] = 3.94
$data[77,4] = 3.83
$data[77,5] = 3.59
$data[77,6] = 3.76
$data[78,0] = 79
$data[78,1] = 230017
$data[78,2] = 'ADHIKARI A.H.C.S.'
$data[78,3] = 3.9
$data[78,4] = 3.95
$data[78,5] = 3.44
$data[78,6] = 3.75
$data[79,0] = 80
$data[79,1] = 230327
$data[79,2] = 'KAUSHALYA R.G.S.P.'
$data[79,3] = 3.85
$data[79,4] = 3.68
$data[79,5] = 3.77
$data[79,6] = 3.75
$data[80,0] = 81
$data[80,1] = 230473
$data[80,2] = 'PERAMUNAGE D.S.'
$data[80,3] = 3.89
$data[80,4] = 3.92
$data[80,5] = 3.44
$data[80,6] = 3.74
$data[81,0] = 82
$data[81,1] = 230727
$data[81,2] = 'WIJESINGHE W.A.P.W.'
$data[81,3] = 3.79
$data[81,4] = 3.87
$data[81,5] = 3.54
$data[81,6] = 3.74
$data[82,0] = 83
$data[82,1] = 230012
$data[82,2] = 'ABEYWARDHANA T.C.W.'
$data[82,3] = 3.91
$data[82,4] = 3.8
$data[82,5] = 3.56
$data[82,6] = 3.73
$data[83,0] = 84
$data[83,1] = 230147
$data[83,2] = 'DILHARA D.S.'
$data[83,3] = 3.89
$data[83,4] = 3.83
$data[83,5] = 3.51
$data[83,6] = 3.73
$data[84,0] = 85
$data[84,1] = 230077
$data[84,2] = 'BANDARA K.M.N.D.'
$data[84,3] = 3.79
$data[84,4] = 3.75
$data[84,5] = 3.62
$data[84,6] = 3.72
$data[85,0] = 86
$data[85,1] = 230520
$data[85,2] = 'RANASINGHE A.G.N.S.'
$data[85,3] = 3.85
$data[85,4] = 3.85
$data[85,5] = 3.49
$data[85,6] = 3.72
$data[86,0] = 87
$data[86,1] = 230495
$data[86,2] = 'PRABHARSHA H.W.D.'
$data[86,3] = 3.85
$data[86,4] = 3.92
$data[86,5] = 3.42
$data[86,6] = 3.71
$data[87,0] = 88
$data[87,1] = 230375
$data[87,2] = 'LENMINI B.L.W.'
$data[87,3] = 3.85
$data[87,4] = 3.71
$data[87,5] = 3.6
$data[87,6] = 3.7
$data[88,0] = 89
$data[88,1] = 230407
$data[88,2] = 'MEEDENIYA M.M.H.'
$data[88,3] = 4
$data[88,4] = 3.75
$data[88,5] = 3.47
$data[88,6] = 3.7
$data[89,0] = 90
$data[89,1] = 230444
$data[89,2] = 'NIRMANI W.T.'
$data[89,3] = 3.79
$data[89,4] = 3.62
$data[89,5] = 3.76
$data[89,6] = 3.7
$data[90,0] = 91
$data[90,1] = 230261
$data[90,2] = 'INDUWARA M.L.A.S.'
$data[90,3] = 4
$data[90,4] = 3.77
$data[90,5] = 3.42
$data[90,6] = 3.69
$data[91,0] = 92
$data[91,1] = 230458
$data[91,2] = 'PALIHENA H.H.'
$data[91,3] = 3.96
$data[91,4] = 3.74
$data[91,5] = 3.47
$data[91,6] = 3.69
$data[92,0] = 93
$data[92,1] = 230527
$data[92,2] = 'RANAWAKA R.A.G.K.'
$data[92,3] = 4
$data[92,4] = 3.86
$data[92,5] = 3.34
$data[92,6] = 3.69
$data[93,0] = 94
$data[93,1] = 230248
$data[93,2] = 'HIMASARA W.V.M.J.'
$data[93,3] = 3.94
$data[93,4] = 3.66
$data[93,5] = 3.54
$data[93,6] = 3.68
$data[94,0] = 95
$data[94,1] = 230735
$data[94,2] = 'WITHANAGE G.D.N.'
$data[94,3] = 3.94
$data[94,4] = 3.82
$data[94,5] = 3.4
$data[94,6] = 3.68
$data[95,0] = 96
$data[95,1] = 230229
$data[95,2] = 'HANSINDU M.M.A.D.'
$data[95,3] = 3.85
$data[95,4] = 3.78
$data[95,5] = 3.4
$data[95,6] = 3.67
$data[96,0] = 97
$data[96,1] = 230581
$data[96,2] = 'SANTHOSH S.'
$data[96,3] = 3.79
$data[96,4] = 3.68
$data[96,5] = 3.55
$data[96,6] = 3.67
$data[97,0] = 98
$data[97,1] = 230650
$data[97,2] = 'UBEYSEKARA V.T.T.'
$data[97,3] = 4
$data[97,4] = 3.74
$data[97,5] = 3.37
$data[97,6] = 3.66
$data[98,0] = 99
$data[98,1] = 230013
$data[98,2] = 'ABEYWARNA D.H.'
$data[98,3] = 3.85
$data[98,4] = 3.6
$data[98,5] = 3.59
$data[98,6] = 3.65
$data[99,0] = 100
$data[99,1] = 230175
$data[99,2] = 'ERANGA W.A.O.'
$data[99,3] = 3.96
$data[99,4] = 3.7
$data[99,5] = 3.4
$data[99,6] = 3.65
$data[100,0] = 101
$data[100,1] = 230208
$data[100,2] = 'GUNASEKARA H.M.'
$data[100,3] = 4
$data[100,4] = 3.64
$data[100,5] = 3.46
$data[100,6] = 3.65
$data[101,0] = 102
$data[101,1] = 230238
$data[101,2] = 'HENDENIYA H.M.J.C.'
$data[101,3] = 3.76
$data[101,4] = 3.48
$data[101,5] = 3.75
$data[101,6] = 3.63
$data[102,0] = 103
$data[102,1] = 230563
$data[102,2] = 'SAMARANAYAKA H.D.J.D.'
$data[102,3] = 3.89
$data[102,4] = 3.81
$data[102,5] = 3.3
$data[102,6] = 3.63
$data[103,0] = 104
$data[103,1] = 230493
$data[103,2] = 'PIYUMAL N.P.P.'
$data[103,3] = 3.94
$data[103,4] = 3.69
$data[103,5] = 3.33
$data[103,6] = 3.61
$data[104,0] = 105
$data[104,1] = 230259
$data[104,2] = 'IMBULPITIYA B.N.'
$data[104,3] = 3.85
$data[104,4] = 3.57
$data[104,5] = 3.45
$data[104,6] = 3.59
$data[105,0] = 106
$data[105,1] = 230395
$data[105,2] = 'MANATUNGA K.D.'
$data[105,3] = 3.85
$data[105,4] = 3.68
$data[105,5] = 3.33
$data[105,6] = 3.59
$data[106,0] = 107
$data[106,1] = 230730
$data[106,2] = 'WIJETHILAKA J.S.'
$data[106,3] = 4
$data[106,4] = 3.65
$data[106,5] = 3.17
$data[106,6] = 3.59
$data[107,0] = 108
$data[107,1] = 230033
$data[107,2] = 'AMARASINGHE A.A.D.K.'
$data[107,3] = 3.85
$data[107,4] = 3.64
$data[107,5] = 3.35
$data[107,6] = 3.58
$data[108,0] = 109
$data[108,1] = 230636
$data[108,2] = 'THARUSHIKA G.K.E.'
$data[108,3] = 3.89
$data[108,4] = 3.74
$data[108,5] = 3.24
$data[108,6] = 3.58
$data[109,0] = 110
$data[109,1] = 230183
$data[109,2] = 'FERNANDO LTJ'
$data[109,3] = 4
$data[109,4] = 3.54
$data[109,5] = 3.25
$data[109,6] = 3.57
$data[110,0] = 111
$data[110,1] = 230564
$data[110,2] = 'SAMARASEKARA S.M.R.P.'
$data[110,3] = 3.87
$data[110,4] = 3.57
$data[110,5] = 3.22
$data[110,6] = 3.51
$data[111,0] = 112
$data[111,1] = 230203
$data[111,2] = 'GUNARATHNA K.T.M.B.'
$data[111,3] = 3.85
$data[111,4] = 3.49
$data[111,5] = 3.25
$data[111,6] = 3.5
$data[112,0] = 113
$data[112,1] = 230224
$data[112,2] = 'HAKAM M.R.A.'
$data[112,3] = 3.85
$data[112,4] = 3.53
$data[112,5] = 3.22
$data[112,6] = 3.49
$data[113,0] = 114
$data[113,1] = 230268
$data[113,2] = 'JAYAKODY J.A.C.P.'
$data[113,3] = 3.85
$data[113,4] = 3.62
$data[113,5] = 3.1
$data[113,6] = 3.48
$data[114,0] = 115
$data[114,1] = 230449
$data[114,2] = 'NUWANAKA W.A.S.'
$data[114,3] = 3.81
$data[114,4] = 3.32
$data[114,5] = 2.83
$data[114,6] = 3.27

$ws.Range("A2:G116").Value = $data

# --- Column C width was narrowed slightly ---
$ws.Columns.Item(3).ColumnWidth = 16.26

# --- Selection / scroll position left where the last edit was made ---
$ws.Range("B115").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 111
